$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 5 and row 6 (columns A:E), simulating a sort
# that reorders these two data rows.
$row5 = $ws.Range("A5:E5").Value2
$row6 = $ws.Range("A6:E6").Value2

$ws.Range("A5:E5").Value2 = $row6
$ws.Range("A6:E6").Value2 = $row5
